$d = $word.ActiveDocument

# --- Step 1: Insert a new paragraph "RN 004.9. Um módulo pode ser constituído de
#     Sistema Embarcado." right before the existing paragraph that reads
#     "RN 004.9. Um módulo pode ser constituído de API, Banco de Dados e Sistema
#     Embarcado." ---
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "RN 004.9.*Banco de Dados*") {
        $targetIdx = $i
        break
    }
}

$insertionRange = $d.Paragraphs.Item($targetIdx).Range
$insertionRange.Collapse(1) # wdCollapseStart
$insertionRange.InsertParagraphBefore()

# The freshly inserted (still empty) paragraph now occupies the old index;
# everything else shifted down by one.
$d.Paragraphs.Item($targetIdx).Range.Text = "RN 004.9. Um módulo pode ser constituído de Sistema Embarcado."

# --- Step 2: Renumber the (now shifted) "RN 004.9" paragraph (API, Banco de
#     Dados e Sistema Embarcado) to "RN 004.10". Edit surgically - replace just
#     the lone "9" character/run with "10" - instead of rewriting the whole
#     paragraph, so the untouched "RN", " 00", "4" runs keep their identity. ---
$renumIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "RN 004.9.*Banco de Dados*") {
        $renumIdx = $i
        break
    }
}
$renumPara = $d.Paragraphs.Item($renumIdx)
$nineStart = $renumPara.Range.Start + 7   # length of "RN 004." == 7
$nineRange = $d.Range($nineStart, $nineStart + 1)
$nineRange.Text = "10"

# --- Step 3: Delete the old "RN 004.10. Um módulo pode ser constituído de API e
#     Sistema Embarcado." paragraph entirely (now redundant/duplicate numbering). ---
$delIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "RN 004.10.*API e Sistema Embarcado*") {
        $delIdx = $i
        break
    }
}
$d.Paragraphs.Item($delIdx).Range.Delete()
